$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.037.28'
$ws.Range("E2").Value = '  -3.08%  '
$ws.Range("D3").Value = '1.798.89'
$ws.Range("E3").Value = '  -3.27%  '
$ws.Range("E4").Value = '  -0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '307.23'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.13%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4193'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -3.11%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3578'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -3.89%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07090'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -3.94%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8470'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -4.17%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.18'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -5.22%  '
$ws.Range("D12").Value = '1.801.52'
$ws.Range("E12").Value = '  -4.06%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.299'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -3.93%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.355'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.35%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.06753'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.13%  '
$ws.Range("E16").Value = '  +0.11%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '79.98'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.16%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008698'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -4.55%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.11%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '15.01'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.98%  '
$ws.Range("D21").Value = '26.801.95'
$ws.Range("E21").Value = '  -4.31%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.056'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.10%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.98'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("D24").Value = '1.965.37'
$ws.Range("E24").Value = '  -6.25%  '
$ws.Range("E25").Value = '  -2.74%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '152.83'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.50%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.11'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -5.26%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '5.015'
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '112.93'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -3.07%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.647'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -12.30%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.09015'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.47%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.7214'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -9.18%  '
$ws.Range("E33").Value = '  -4.47%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.302'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -7.58%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.083'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -8.56%  '
$ws.Range("E36").Value = '  -0.02%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.080'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.88%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01903'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.59%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05123'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -6.47%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.1626'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -4.16%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.4954'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -5.05%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.576'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -9.58%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '8.052'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -8.16%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.968'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -12.33%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '104.81'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.26%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.16'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -4.80%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.06293'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.47%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.4530'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -6.00%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.603'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.61%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.697'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -9.83%  '
